function Test-ExactMatch($str, $target) {
    if ($str.Length -ne $target.Length) { return $false }
    $sc = $str.ToCharArray()
    $tc = $target.ToCharArray()
    for ($i = 0; $i -lt $sc.Length; $i++) {
        if ([int]$sc[$i] -ne [int]$tc[$i]) { return $false }
    }
    return $true
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$rowCount = $usedRange.Rows.Count

for ($r = 1; $r -le $rowCount; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2
    if ($val -ne $null -and $val -ne "") {
        $parts = $val -split ',\s*'
        $nonSystem = @()
        $systemCount = 0
        foreach ($p in $parts) {
            if (Test-ExactMatch $p "System") {
                $systemCount++
            } else {
                $nonSystem += $p
            }
        }
        if ($systemCount -gt 0) {
            $newParts = $nonSystem
            for ($i = 0; $i -lt $systemCount; $i++) {
                $newParts += "System"
            }
            $newVal = $newParts -join ", "
            if ($newVal -ne $val) {
                $cell.Value2 = $newVal
            }
        }
    }
}
